# PRJ0019297 CVAS Time Tracking Changes+ changes after SB Referesh
$wb = $excel.ActiveWorkbook

$usersSheet = $wb.Worksheets.Item("Users")

# Update the Global Search User value from "Coartney Williams" to "Coartney Trone"
$usersSheet.Range("A2").Value = "Coartney Trone"

# Bold the header cell
$usersSheet.Range("A1").Font.Bold = $true

# Make the Users sheet the active/selected sheet
$usersSheet.Activate() | Out-Null
$usersSheet.Range("A1").Select() | Out-Null
